# Auto-generated script to update FFXIV Goblin Profits market-data values
# Applies per-cell numeric updates across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6229.4253
$ws.Range("I64").Value = 3159.2334
$ws.Range("J64").Value = 11647.412
$ws.Range("K64").Value = 3159.2334
$ws.Range("L64").Value = 11647.412
$ws.Range("M64").Value = -2911.2334
$ws.Range("N64").Value = -12143.412

$ws.Range("H67").Value = 6229.4253
$ws.Range("I67").Value = 3159.2334
$ws.Range("J67").Value = 11647.412
$ws.Range("K67").Value = 3159.2334
$ws.Range("L67").Value = 11647.412
$ws.Range("M67").Value = -2301.2334
$ws.Range("N67").Value = -13363.412

$ws.Range("H74").Value = 3763
$ws.Range("I74").Value = 3499.25
$ws.Range("J74").Value = 4466.3335
$ws.Range("K74").Value = 3499.25
$ws.Range("L74").Value = 4466.3335
$ws.Range("M74").Value = -2563.25
$ws.Range("N74").Value = -6338.3335

$ws.Range("H77").Value = 3763
$ws.Range("I77").Value = 3499.25
$ws.Range("J77").Value = 4466.3335
$ws.Range("K77").Value = 17496.25
$ws.Range("L77").Value = 22331.6675
$ws.Range("M77").Value = -12816.25
$ws.Range("N77").Value = -31691.6675

$ws.Range("H107").Value = 506713.03
$ws.Range("I107").Value = 655352.25
$ws.Range("J107").Value = 1339.8
$ws.Range("K107").Value = 655352.25
$ws.Range("L107").Value = 1339.8
$ws.Range("M107").Value = -653432.25
$ws.Range("N107").Value = -5179.8

$ws.Range("H132").Value = 1436.8572
$ws.Range("I132").Value = 1422.9048
$ws.Range("J132").Value = 1520.5714
$ws.Range("K132").Value = 4268.7144
$ws.Range("L132").Value = 4561.7142
$ws.Range("M132").Value = -1738.7144
$ws.Range("N132").Value = -9621.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6728.5
$ws.Range("I32").Value = 4081.4524
$ws.Range("K32").Value = 4081.4524
$ws.Range("M32").Value = -3794.4524

$ws.Range("H36").Value = 2704
$ws.Range("I36").Value = 2704
$ws.Range("K36").Value = 2704
$ws.Range("M36").Value = -2358

$ws.Range("H61").Value = 5410.4644
$ws.Range("J61").Value = 9949.5
$ws.Range("L61").Value = 9949.5
$ws.Range("N61").Value = -10373.5

$ws.Range("H63").Value = 4421.0625
$ws.Range("I63").Value = 1767.2727
$ws.Range("J63").Value = 10259.4
$ws.Range("K63").Value = 1767.2727
$ws.Range("L63").Value = 10259.4
$ws.Range("M63").Value = -1081.2727
$ws.Range("N63").Value = -11631.4

$ws.Range("H66").Value = 4421.0625
$ws.Range("I66").Value = 1767.2727
$ws.Range("J66").Value = 10259.4
$ws.Range("K66").Value = 8836.363499999999
$ws.Range("L66").Value = 51297
$ws.Range("M66").Value = -5404.363499999999
$ws.Range("N66").Value = -58161

$ws.Range("H102").Value = 4575.7915
$ws.Range("I102").Value = 1321.3334
$ws.Range("K102").Value = 1321.3334
$ws.Range("M102").Value = 300.6666

$ws.Range("H136").Value = 5410.4644
$ws.Range("J136").Value = 9949.5
$ws.Range("L136").Value = 29848.5
$ws.Range("N136").Value = -34948.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5966.074
$ws.Range("I20").Value = 7996.9443
$ws.Range("J20").Value = 1904.3334
$ws.Range("K20").Value = 7996.9443
$ws.Range("L20").Value = 1904.3334
$ws.Range("M20").Value = -7749.9443
$ws.Range("N20").Value = -2398.3334

$ws.Range("H82").Value = 36752.9
$ws.Range("I82").Value = 22921.5
$ws.Range("J82").Value = 57500
$ws.Range("K82").Value = 22921.5
$ws.Range("L82").Value = 57500
$ws.Range("M82").Value = -22538.5
$ws.Range("N82").Value = -58266

$ws.Range("H85").Value = 36752.9
$ws.Range("I85").Value = 22921.5
$ws.Range("J85").Value = 57500
$ws.Range("K85").Value = 22921.5
$ws.Range("L85").Value = 57500
$ws.Range("M85").Value = -21595.5
$ws.Range("N85").Value = -60152

$ws.Range("H105").Value = 3913.361
$ws.Range("I105").Value = 3859.697
$ws.Range("K105").Value = 3859.697
$ws.Range("M105").Value = -2112.697

$ws.Range("H134").Value = 2296.6924
$ws.Range("I134").Value = 2088.2188
$ws.Range("K134").Value = 6264.6564
$ws.Range("M134").Value = -3729.6564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3329.2354
$ws.Range("I58").Value = 3649.9167
$ws.Range("K58").Value = 3649.9167
$ws.Range("M58").Value = -3446.9167

$ws.Range("H69").Value = 44915
$ws.Range("I69").Value = 24830
$ws.Range("K69").Value = 24830
$ws.Range("M69").Value = -24081

$ws.Range("H72").Value = 44915
$ws.Range("I72").Value = 24830
$ws.Range("K72").Value = 74490
$ws.Range("M72").Value = -70746

$ws.Range("H136").Value = 3329.2354
$ws.Range("I136").Value = 3649.9167
$ws.Range("K136").Value = 10949.7501
$ws.Range("M136").Value = -8399.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 166.71428
$ws.Range("J12").Value = 212.16667
$ws.Range("L12").Value = 636.50001
$ws.Range("N12").Value = -982.50001

$ws.Range("H108").Value = 6648.9
$ws.Range("I108").Value = 761.125
$ws.Range("J108").Value = 30200
$ws.Range("K108").Value = 2283.375
$ws.Range("L108").Value = 90600
$ws.Range("M108").Value = 596.625
$ws.Range("N108").Value = -96360

$ws.Range("H113").Value = 1362.8572
$ws.Range("J113").Value = 747.5
$ws.Range("L113").Value = 2242.5
$ws.Range("N113").Value = -6582.5

$ws.Range("H117").Value = 2084.5833
$ws.Range("I117").Value = 2725.6
$ws.Range("J117").Value = 1626.7142
$ws.Range("K117").Value = 8176.799999999999
$ws.Range("L117").Value = 4880.142599999999
$ws.Range("M117").Value = -4734.799999999999
$ws.Range("N117").Value = -11764.1426

$ws.Range("H121").Value = 1684.7931
$ws.Range("I121").Value = 454.63635
$ws.Range("K121").Value = 1363.90905
$ws.Range("M121").Value = -53.90904999999998

$ws.Range("H137").Value = 13524.125
$ws.Range("I137").Value = 9565.5
$ws.Range("J137").Value = 15899.3
$ws.Range("K137").Value = 28696.5
$ws.Range("L137").Value = 47697.89999999999
$ws.Range("M137").Value = -23596.5
$ws.Range("N137").Value = -57897.89999999999

$ws.Range("H138").Value = 38464170
$ws.Range("I138").Value = 50001320
$ws.Range("J138").Value = 6998.6665
$ws.Range("K138").Value = 150003960
$ws.Range("L138").Value = 20995.9995
$ws.Range("M138").Value = -149998820
$ws.Range("N138").Value = -31275.9995

$ws.Range("H139").Value = 5733.706
$ws.Range("I139").Value = 7578.8335
$ws.Range("K139").Value = 22736.5005
$ws.Range("M139").Value = -17596.5005

$ws.Range("H141").Value = 125005140
$ws.Range("I141").Value = 250003340
$ws.Range("K141").Value = 750010020
$ws.Range("M141").Value = -750004840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1526.1923
$ws.Range("I97").Value = 1429.8695
$ws.Range("K97").Value = 1429.8695
$ws.Range("M97").Value = -933.8695

$ws.Range("H107").Value = 1286.625
$ws.Range("I107").Value = 470.4375
$ws.Range("J107").Value = 2102.8125
$ws.Range("K107").Value = 470.4375
$ws.Range("L107").Value = 2102.8125
$ws.Range("M107").Value = 1449.5625
$ws.Range("N107").Value = -5942.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 8120.9165
$ws.Range("I100").Value = 4510
$ws.Range("J100").Value = 10700.143
$ws.Range("K100").Value = 4510
$ws.Range("L100").Value = 10700.143
$ws.Range("M100").Value = -3969
$ws.Range("N100").Value = -11782.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 3500
$ws.Range("J9").Value = 3500
$ws.Range("L9").Value = 3500
$ws.Range("N9").Value = -3780

$ws.Range("H74").Value = 26208.25
$ws.Range("I74").Value = 22998
$ws.Range("J74").Value = 26666.857
$ws.Range("K74").Value = 22998
$ws.Range("L74").Value = 26666.857
$ws.Range("M74").Value = -22062
$ws.Range("N74").Value = -28538.857

$ws.Range("H77").Value = 26208.25
$ws.Range("I77").Value = 22998
$ws.Range("J77").Value = 26666.857
$ws.Range("K77").Value = 68994
$ws.Range("L77").Value = 80000.571
$ws.Range("M77").Value = -64314
$ws.Range("N77").Value = -89360.571

$ws.Range("H100").Value = 1210.9286
$ws.Range("I100").Value = 573
$ws.Range("J100").Value = 1565.3334
$ws.Range("K100").Value = 1146
$ws.Range("L100").Value = 3130.6668
$ws.Range("M100").Value = -605
$ws.Range("N100").Value = -4212.6668

$ws.Range("H110").Value = 171248.75
$ws.Range("J110").Value = 171248.75
$ws.Range("L110").Value = 171248.75
$ws.Range("N110").Value = -179428.75

$ws.Range("H132").Value = 3289.96
$ws.Range("I132").Value = 3085.238
$ws.Range("K132").Value = 9255.714
$ws.Range("M132").Value = -6725.714
